# Apply the bank-statement update described by the diff:
#  - SALDO AWAL (opening balance) changes
#  - Existing transaction rows 3-7 get new dates/descriptions/amounts
#  - Row 4's Kategori/Pecah columns become blank
#  - Four new transaction rows are inserted before the closing SALDO AKHIR row
#  - The closing SALDO AKHIR row moves from row 8 to row 12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text (e.g. "01-Jan-2025") is stored as literal text,
# not auto-converted to a date serial by Excel.
$dateCells = @("A3", "A4", "A5", "A6", "A7", "A8", "A9", "A10", "A11")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 2: opening balance (SALDO AWAL) ---
$ws.Range("G2").Value = 496264489

# --- Row 3 ---
$ws.Range("A3").Value = "01-Jan-2025"
$ws.Range("B3").Value = "BY ADMINISTRASI"
$ws.Range("C3").Value = "Tanpa Kategori"
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = 10000
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 496254489

# --- Row 4 ---
$ws.Range("A4").Value = "02-Jan-2025"
$ws.Range("B4").Value = "TRANSFER KE SIMSEM  Tanpa Kategori"
$ws.Range("C4:D4").ClearContents()
$ws.Range("E4").Value = 20138
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 496234351

# --- Row 5 ---
$ws.Range("A5").Value = "03-Jan-2025"
$ws.Range("B5").Value = "KREDIT LAIN-LAIN 050"
$ws.Range("C5").Value = "Tanpa Kategori"
$ws.Range("D5").Value = "No"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 157751823
$ws.Range("G5").Value = 653986174

# --- Row 6 ---
$ws.Range("A6").Value = "10-Jan-2025"
$ws.Range("B6").Value = "RTGS KE NI KETUT"
$ws.Range("C6").Value = "Tanpa Kategori"
$ws.Range("D6").Value = "No"
$ws.Range("E6").Value = 150000000
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 503986174

# --- Row 7 ---
$ws.Range("A7").Value = "10-Jan-2025"
$ws.Range("B7").Value = "RTGS KE NI KETUT"
$ws.Range("C7").Value = "Tanpa Kategori"
$ws.Range("D7").Value = "No"
$ws.Range("E7").Value = 25000
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 503961174

# --- Row 8 (new) ---
$ws.Range("A8").Value = "31-Jan-2025"
$ws.Range("B8").Value = "JASA GIROIBUNGA"
$ws.Range("C8").Value = "Tanpa Kategori"
$ws.Range("D8").Value = "No"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 431953
$ws.Range("G8").Value = 504393127

# --- Row 9 (new) ---
$ws.Range("A9").Value = "31-Jan-2025"
$ws.Range("B9").Value = "PPH"
$ws.Range("C9").Value = "Tanpa Kategori"
$ws.Range("D9").Value = "No"
$ws.Range("E9").Value = 86391
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 504306736

# --- Row 10 (new) ---
$ws.Range("A10").Value = "31-Jan-2025"
$ws.Range("B10").Value = "BY ADMINISTRASI"
$ws.Range("C10").Value = "Tanpa Kategori"
$ws.Range("D10").Value = "No"
$ws.Range("E10").Value = 12000
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 504294736

# --- Row 11 (was row 7's BIAYA ADM KARTU transaction, now re-dated) ---
$ws.Range("A11").Value = "31-Jan-2025"
$ws.Range("B11").Value = "BIAYA ADM KARTU"
$ws.Range("C11").Value = "Tanpa Kategori"
$ws.Range("D11").Value = "No"
$ws.Range("E11").Value = 10000
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 504284736

# --- Row 12 (closing SALDO AKHIR row, moved down from row 8) ---
$ws.Range("A12:D12").ClearContents()
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 504284736
$ws.Range("H12").Value = "SALDO AKHIR"

# Clear the stale SALDO AKHIR text that used to live on row 8 (now a normal
# transaction row) and reset any leftover quote-prefix styling so freshly
# typed text cells keep the workbook's default (unstyled) formatting.
$ws.Range("H8").ClearContents()
foreach ($addr in $dateCells) {
    $ws.Range($addr).Style = "Normal"
}
